$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''256.78'
$ws.Range("E2").Value = '''-0.64%'
$ws.Range("E3").Value = '''0.42%'
$ws.Range("D4").Value = '''4.389'
$ws.Range("E4").Value = '''-8.57%'
$ws.Range("D5").Value = '''0.05889'
$ws.Range("E5").Value = '''-1.35%'
$ws.Range("D6").Value = '''6.629'
$ws.Range("D7").Value = '''0.8547'
$ws.Range("D8").Value = '''0.9371'
$ws.Range("E8").Value = '''-1.55%'
$ws.Range("E9").Value = '''-1.87%'
$ws.Range("D10").Value = '''0.04650'
$ws.Range("E10").Value = '''28.74%'
$ws.Range("D11").Value = '''0.07075'
$ws.Range("E11").Value = '''-1.52%'
$ws.Range("D12").Value = '''0.03072'
$ws.Range("E12").Value = '''-2.87%'
$ws.Range("D13").Value = '''0.09114'
$ws.Range("E13").Value = '''-1.33%'
$ws.Range("D14").Value = '''0.001538'
$ws.Range("E14").Value = '''-0.58%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006159'
$ws.Range("E15").Value = '''3.01%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.493'
$ws.Range("E16").Value = '''0.21%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.186'
$ws.Range("E17").Value = '''-1.31%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.226'
$ws.Range("E18").Value = '''0.30%'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '''0.0006031'
$ws.Range("E19").Value = '''-0.62%'
$ws.Range("E20").Value = '''-2.70%'
$ws.Range("D21").Value = '''0.1269'
$ws.Range("E21").Value = '''-1.59%'
$ws.Range("D22").Value = '''3.910'
$ws.Range("E22").Value = '''10.81%'
$ws.Range("D23").Value = '''0.04279'
$ws.Range("E23").Value = '''1.35%'
$ws.Range("D24").Value = '''0.001219'
$ws.Range("E24").Value = '''0.00%'
$ws.Range("D25").Value = '''0.004284'
$ws.Range("E25").Value = '''-5.08%'
$ws.Range("E26").Value = '''0.03%'
$ws.Range("D27").Value = '''0.0001523'
$ws.Range("E27").Value = '''2.07%'
$ws.Range("D40").Value = '''0.03820'
$ws.Range("E40").Value = '''-0.65%'
$ws.Range("D41").Value = '''0.006240'
$ws.Range("E41").Value = '''54.64%'
$ws.Range("E42").Value = '''-0.26%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002200'
$ws.Range("E43").Value = '''-4.32%'
$ws.Range("B44").Value = 'LocalTraders'
$ws.Range("C44").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D44").Value = '''0.01409'
$ws.Range("E44").Value = '''27.13%'
$ws.Range("D45").Value = '''0.00005359'
$ws.Range("E45").Value = '''-2.42%'
$ws.Range("E46").Value = '''0.03%'
$ws.Range("D47").Value = '''0.06588'
$ws.Range("E48").Value = '''11,741.09%'
$ws.Range("E49").Value = '''0.03%'
$ws.Range("E50").Value = '''0.03%'
